# -----------------------------------------------------------------
# DoE_TempFlume.xlsx -- "Debuggin Parallel script, problems with
# restart threads"
#
#  * Parameters sheet: add a "thickness in" column and refresh the
#    metric thickness conversions for both scale blocks.
#  * Tests Scale 1-21 sheet: record completion/duration for the tests
#    that finished while debugging the restart-thread issue, and fix
#    up the duration that had been logged incorrectly for test 44.
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "Parameters" sheet
# ---------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("Parameters")

# New header cell I1, formatted like the rest of the header row
# (bold / centered / bordered) by copying H1's format over first.
$wsParams.Range("H1").Copy()
$wsParams.Range("I1").PasteSpecial(-4122)
$wsParams.Range("I1").Value = "Unnamed: 7"

# New "thickness" unit column, same layout as the ft/mts columns.
$wsParams.Range("I2").Value = "thickness in"
$wsParams.Range("I3").Value = 0.125

# Recomputed thickness-in-meters for the Scale 1/21 block.
$wsParams.Range("C14").Value = 0.0050038
$wsParams.Range("D14").Value = 0.0074422
$wsParams.Range("E14").Value = 0.0175006

# Recomputed thickness-in-meters for the Scale 1/25 block.
$wsParams.Range("C28").Value = 0.004699
$wsParams.Range("D28").Value = 0.006223
$wsParams.Range("E28").Value = 0.0165862

# ---------------------------------------------------------------
# "Tests Scale 1-21" sheet
# ---------------------------------------------------------------
$wsTests21 = $wb.Worksheets.Item("Tests Scale 1-21")

# Rows that finished running, with their elapsed time (minutes).
$finishedRuns = @{
    4  = 1
    5  = 1
    6  = 1
    7  = 0.5
    8  = 0.5
    9  = 0.5
    10 = 0.5
    11 = 0.5
    12 = 0.5
    13 = 0.2
}

foreach ($row in $finishedRuns.Keys) {
    $wsTests21.Range("E$row").Value = "Completed"
    $wsTests21.Range("F$row").Value = $finishedRuns[$row]
}

# Corrected duration recorded for test 44 (row 46).
$wsTests21.Range("F46").Value = 30
